$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly-run calibration scenarios for years 2014 and 2015 (rows 6 and 7)
$ws.Range("B6:D6").Value = 1
$ws.Range("B7:D7").Value = 1

# Update the active selection to reflect where the user ended up (H13)
$ws.Range("H13").Select()
